$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.841771125793457
$ws.Range("B1").Value = 1.287475824356079
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.548931241035461
